# Topic 1 - Question 1: center-align the "Título da Monografia" paragraph
# on the title slide (slide 1, shape 1 "Retângulo 1").
#
# The shape's text box has three paragraphs:
#   1) "Apresentação Grupo X"                              (already centered)
#   2) (empty paragraph, just an endParaRPr)
#   3) "Título da Monografia: <...>"                        (needs centering)
#
# This adds <a:pPr algn="ctr"/> to paragraph 3's pPr.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tr  = $shp.TextFrame.TextRange

$ppAlignCenter = 2
$titlePara = $tr.Paragraphs(3, 1)
$titlePara.ParagraphFormat.Alignment = $ppAlignCenter
